$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$officesText = "2.1% CR/LFM+CDL/H:2/Offices`n12.0% CR/LFM+CDL/HBET:3-5/Offices`n2.1% CR/LFM+CDN/H:2/Offices`n12.0% CR/LFM+CDN/HBET:3-5/Offices`n6.5% MR/LWAL+CDL/H:1/Offices`n12.9% MR/LWAL+CDL/H:2/Offices`n45.2% MR/LWAL+CDL/HBET:3-5/Offices`n5.2% S/LFM+CDL/HBET:3-5/Offices`n2.0% S/LWAL+CDL/HBET:3-5/Offices`n0.0% CR/LFM+CDL/H:1/Offices`n0.0% CR/LFM+CDN/H:1/Offices`n0.0% S/LFM+CDL/H:1/Offices`n0.0% S/LFM+CDL/H:2/Offices`n0.0% S/LWAL+CDL/H:1/Offices`n0.0% S/LWAL+CDL/H:2/Offices`n0.0% W/LFM+CDL/H:1/Offices`n0.0% W/LFM+CDL/H:2/Offices"

$tradeText = "1.2% CR/LFM+CDL/H:2/Trade`n0.0% CR/LFM+CDN/HBET:3-5/Trade`n1.2% CR/LFM+CDL/H:2/Trade`n0.0% CR/LFM+CDN/HBET:3-5/Trade`n58.2% MR/LWAL+CDL/H:1/Trade`n6.5% MR/LWAL+CDL/H:2/Trade`n0.0% MR/LWAL+CDL/HBET:3-5/Trade`n0.0% S/LFM+CDL/HBET:3-5/Trade`n0.0% S/LWAL+CDL/HBET:3-5/Trade`n10.8% CR/LFM+CDL/H:1/Trade`n10.8% CR/LFM+CDN/H:1/Trade`n4.6% S/LFM+CDL/H:1/Trade`n0.5% S/LFM+CDL/H:2/Trade`n1.8% S/LWAL+CDL/H:1/Trade`n0.2% S/LWAL+CDL/H:2/Trade`n3.8% W/LFM+CDL/H:1/Trade`n0.4% W/LFM+CDL/H:2/Trade"

$hotelsText = "1.2% CR/LFM+CDL/H:2/Hotels`n8.4% CR/LFM+CDL/HBET:3-5/Hotels`n1.2% CR/LFM+CDN/H:2/Hotels`n8.4% CR/LFM+CDN/HBET:3-5/Hotels`n12.9% MR/LWAL+CDL/H:1/Hotels`n 6.5% MR/LWAL+CDL/H:2/Hotels`n 45.3% MR/LWAL+CDL/HBET:3-5/Hotels`n 5.2% S/LFM+CDL/HBET:3-5/Hotels`n 2.0% S/LWAL+CDL/HBET:3-5/Hotels`n 2.4% CR/LFM+CDL/H:1/Hotels`n 2.4% CR/LFM+CDN/H:1/Hotels`n 0.0% S/LFM+CDL/H:1/Hotels`n 0.0% S/LFM+CDL/H:2/Hotels`n 0.0% S/LWAL+CDL/H:1/Hotels`n 0.0% S/LWAL+CDL/H:2/Hotels`n 0.8% W/LFM+CDL/H:1/Hotels`n 3.3% W/LFM+CDL/H:2/Hotels"

# Update cell contents: B2 = Offices (mix of CDL/CDN), C2 = Trade, D2 = Hotels
$ws.Range("B2").Value = $officesText
$ws.Range("C2").Value = $tradeText
$ws.Range("D2").Value = $hotelsText

# Apply wrap text style to the data row
$ws.Range("B2:D2").WrapText = $true

# Set column widths (values chosen so stored width matches target after Excel's internal padding)
$ws.Range("B1").EntireColumn.ColumnWidth = 35.666666666666664
$ws.Range("C1").EntireColumn.ColumnWidth = 32.166666666666664
$ws.Range("D1").EntireColumn.ColumnWidth = 34.498697916666664
$ws.Range("E1:G1").EntireColumn.ColumnWidth = 50.830729166666664

# Set row height for the data row to match auto-fit height with wrapped text
$ws.Rows.Item(2).RowHeight = 272

# Match the active selection left by the original author
$ws.Range("B2").Select() | Out-Null

Write-Output "done"
